$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E8").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G8").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H8").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K8").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L8").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E11").Value = 13128994
$ws.Range("F11").Value = 4499875
$ws.Range("G11").Value = 4432450
$ws.Range("H11").Value = -4433
$ws.Range("I11").Value = 15408661
$ws.Range("J11").Value = "-"
$ws.Range("K11").Value = 7112
$ws.Range("L11").Value = 7214
$ws.Range("M11").Value = 5031
$ws.Range("N11").Value = 5751
$ws.Range("E12").Value = 13128994
$ws.Range("F12").Value = 4499875
$ws.Range("G12").Value = 4432450
$ws.Range("H12").Value = -4433
$ws.Range("I12").Value = 15408661
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 7112
$ws.Range("L12").Value = 7214
$ws.Range("M12").Value = 5031
$ws.Range("N12").Value = 5751
$ws.Range("E14").Value = 6409330
$ws.Range("F14").Value = 6439655
$ws.Range("G14").Value = 5161840
$ws.Range("H14").Value = 10260688
$ws.Range("I14").Value = 15870434
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = 9273
$ws.Range("L14").Value = 4687
$ws.Range("M14").Value = 6540
$ws.Range("N14").Value = 14797
$ws.Range("E15").Value = 6409330
$ws.Range("F15").Value = 6439655
$ws.Range("G15").Value = 5161840
$ws.Range("H15").Value = 10260688
$ws.Range("I15").Value = 15870434
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 9273
$ws.Range("L15").Value = 4687
$ws.Range("M15").Value = 6540
$ws.Range("N15").Value = 14797
$ws.Range("G18").Value = "-"
$ws.Range("E19").Value = 2750
$ws.Range("F19").Value = 5000
$ws.Range("G19").Value = "-"
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = "-"
$ws.Range("K19").Value = 0
$ws.Range("N19").Value = -55
$ws.Range("E20").Value = 19541074
$ws.Range("F20").Value = 10944530
$ws.Range("G20").Value = 9594290
$ws.Range("H20").Value = 10256255
$ws.Range("I20").Value = 31279095
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 16385
$ws.Range("L20").Value = 11901
$ws.Range("M20").Value = 11571
$ws.Range("N20").Value = 20493
$ws.Range("E24").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F24").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G24").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H24").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I24").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J24").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K24").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L24").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M24").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N24").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E27").Value = 1413754
$ws.Range("F27").Value = 864133
$ws.Range("G27").Value = 1000933
$ws.Range("H27").Value = -1000933
$ws.Range("I27").Value = 1557100
$ws.Range("J27").Value = "-"
$ws.Range("K27").Value = 2868668
$ws.Range("L27").Value = 3042210
$ws.Range("M27").Value = 1696585
$ws.Range("N27").Value = 1831608
$ws.Range("E28").Value = 1413754
$ws.Range("F28").Value = 864133
$ws.Range("G28").Value = 1000933
$ws.Range("H28").Value = -1000933
$ws.Range("I28").Value = 1557100
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 2868668
$ws.Range("L28").Value = 3042210
$ws.Range("M28").Value = 1696585
$ws.Range("N28").Value = 1831608
$ws.Range("E30").Value = 847981
$ws.Range("F30").Value = 1104492
$ws.Range("G30").Value = 1176056
$ws.Range("H30").Value = 1217745
$ws.Range("I30").Value = 1388773
$ws.Range("J30").Value = "-"
$ws.Range("K30").Value = 2802212
$ws.Range("L30").Value = 1574089
$ws.Range("M30").Value = 2039002
$ws.Range("N30").Value = 4570678
$ws.Range("E31").Value = 847981
$ws.Range("F31").Value = 1104492
$ws.Range("G31").Value = 1176056
$ws.Range("H31").Value = 1217745
$ws.Range("I31").Value = 1388773
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2802212
$ws.Range("L31").Value = 1574089
$ws.Range("M31").Value = 2039002
$ws.Range("N31").Value = 4570678
$ws.Range("G34").Value = "-"
$ws.Range("E35").Value = 525
$ws.Range("G35").Value = "-"
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = "-"
$ws.Range("K35").Value = 0
$ws.Range("N35").Value = -19267
$ws.Range("E36").Value = 2262260
$ws.Range("F36").Value = 1969150
$ws.Range("G36").Value = 2176989
$ws.Range("H36").Value = 216812
$ws.Range("I36").Value = 2945873
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 5670880
$ws.Range("L36").Value = 4616299
$ws.Range("M36").Value = 3735587
$ws.Range("N36").Value = 6383019
$ws.Range("E40").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F40").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G40").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H40").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I40").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J40").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K40").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L40").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M40").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N40").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E43").Value = 135058
$ws.Range("F43").Value = 192035
$ws.Range("G43").Value = 225819
$ws.Range("H43").Value = 247211
$ws.Range("I43").Value = 249967
$ws.Range("J43").Value = -78390
$ws.Range("K43").Value = 403356018
$ws.Range("L43").Value = 421709177
$ws.Range("M43").Value = 337226198
$ws.Range("N43").Value = 318485133
$ws.Range("E45").Value = 132304
$ws.Range("F45").Value = 171514
$ws.Range("G45").Value = 227837
$ws.Range("H45").Value = 238586
$ws.Range("I45").Value = 247339
$ws.Range("J45").Value = -110342
$ws.Range("K45").Value = 302190445
$ws.Range("L45").Value = 335841476
$ws.Range("M45").Value = 311774006
$ws.Range("N45").Value = 308892208
$ws.Range("E49").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F49").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G49").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H49").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I49").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J49").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K49").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L49").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M49").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N49").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E52").Value = -401078
$ws.Range("F52").Value = -658775
$ws.Range("G52").Value = -680853
$ws.Range("H52").Value = 680853
$ws.Range("I52").Value = -2838040
$ws.Range("J52").Value = -1241537
$ws.Range("K52").Value = -1973157
$ws.Range("L52").Value = -2527768
$ws.Range("M52").Value = -1306919
$ws.Range("N52").Value = -1678140
$ws.Range("E53").Value = -401078
$ws.Range("F53").Value = -658775
$ws.Range("G53").Value = -680853
$ws.Range("H53").Value = 680853
$ws.Range("I53").Value = -2838040
$ws.Range("J53").Value = -1241537
$ws.Range("K53").Value = -1973157
$ws.Range("L53").Value = -2527768
$ws.Range("M53").Value = -1306919
$ws.Range("N53").Value = -1678140
$ws.Range("E55").Value = -606001
$ws.Range("F55").Value = -735201
$ws.Range("G55").Value = -790868
$ws.Range("H55").Value = -934688
$ws.Range("I55").Value = -1197590
$ws.Range("J55").Value = -1496396
$ws.Range("K55").Value = -2572488
$ws.Range("L55").Value = -822947
$ws.Range("M55").Value = -1551987
$ws.Range("N55").Value = -3840473
$ws.Range("E56").Value = -606001
$ws.Range("F56").Value = -735201
$ws.Range("G56").Value = -790868
$ws.Range("H56").Value = -934688
$ws.Range("I56").Value = -1197590
$ws.Range("J56").Value = -1496396
$ws.Range("K56").Value = -2572488
$ws.Range("L56").Value = -822947
$ws.Range("M56").Value = -1551987
$ws.Range("N56").Value = -3840473
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = -342
$ws.Range("G59").Value = "-"
$ws.Range("G60").Value = "-"
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = "-"
$ws.Range("J60").Value = 0
$ws.Range("E61").Value = -1007079
$ws.Range("F61").Value = -1394318
$ws.Range("G61").Value = -1471721
$ws.Range("H61").Value = -253835
$ws.Range("I61").Value = -4035630
$ws.Range("J61").Value = -2737933
$ws.Range("K61").Value = -4545645
$ws.Range("L61").Value = -3350715
$ws.Range("M61").Value = -2858906
$ws.Range("N61").Value = -5518613
$ws.Range("E65").Value = "فصل سوم منتهی به 1399/09"
$ws.Range("F65").Value = "فصل چهارم منتهی به 1399/12"
$ws.Range("G65").Value = "فصل اول منتهی به 1400/03"
$ws.Range("H65").Value = "فصل دوم منتهی به 1400/06"
$ws.Range("I65").Value = "فصل سوم منتهی به 1400/09"
$ws.Range("J65").Value = "فصل چهارم منتهی به 1400/12"
$ws.Range("K65").Value = "فصل اول منتهی به 1401/03"
$ws.Range("L65").Value = "فصل دوم منتهی به 1401/06"
$ws.Range("M65").Value = "فصل سوم منتهی به 1401/09"
$ws.Range("N65").Value = "فصل چهارم منتهی به 1401/12"
$ws.Range("E68").Value = 162655
$ws.Range("F68").Value = 205358
$ws.Range("G68").Value = 320080
$ws.Range("H68").Value = -320080
$ws.Range("I68").Value = 1016033
$ws.Range("J68").Value = -34433
$ws.Range("K68").Value = 895511
$ws.Range("L68").Value = 514442
$ws.Range("M68").Value = 389666
$ws.Range("N68").Value = 153468
$ws.Range("E69").Value = 162655
$ws.Range("F69").Value = 205358
$ws.Range("G69").Value = 320080
$ws.Range("H69").Value = -320080
$ws.Range("I69").Value = 1016033
$ws.Range("J69").Value = -34433
$ws.Range("K69").Value = 895511
$ws.Range("L69").Value = 514442
$ws.Range("M69").Value = 389666
$ws.Range("N69").Value = 153468
$ws.Range("E71").Value = 241980
$ws.Range("F71").Value = 369291
$ws.Range("G71").Value = 385188
$ws.Range("H71").Value = 283057
$ws.Range("I71").Value = 191183
$ws.Range("J71").Value = 253328
$ws.Range("K71").Value = 229724
$ws.Range("L71").Value = 751142
$ws.Range("M71").Value = 487015
$ws.Range("N71").Value = 730205
$ws.Range("E72").Value = 241980
$ws.Range("F72").Value = 369291
$ws.Range("G72").Value = 385188
$ws.Range("H72").Value = 283057
$ws.Range("I72").Value = 191183
$ws.Range("J72").Value = 253328
$ws.Range("K72").Value = 229724
$ws.Range("L72").Value = 751142
$ws.Range("M72").Value = 487015
$ws.Range("N72").Value = 730205
$ws.Range("E75").Value = 404635
$ws.Range("F75").Value = 574649
$ws.Range("G75").Value = 705268
$ws.Range("H75").Value = -37023
$ws.Range("I75").Value = 1207216
$ws.Range("J75").Value = 218895
$ws.Range("K75").Value = 1125235
$ws.Range("L75").Value = 1265584
$ws.Range("M75").Value = 876681
$ws.Range("N75").Value = 883673
